# Slide 6, shape "object 19" (the MySQL / MariaDB textbox in the bottom
# right corner) loses its second line of text ("MariaDB"), leaving the
# paragraph mark (and its endParaRPr) in place but with no run. Because
# the shape auto-fits its text, its height also shrinks slightly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(9)

# Sanity check we grabbed the right shape before mutating it.
if ($sh.Name -ne "object 19") {
    throw "Expected shape 'object 19' but found '$($sh.Name)'"
}

$tr = $sh.TextFrame.TextRange

# "MySQL" occupies characters 1-5, character 6 is the paragraph break,
# and "MariaDB" occupies characters 7-13. Clearing just that substring
# removes the run but keeps the (now empty) second paragraph intact.
$mariaDb = $tr.Characters(7, 7)
$mariaDb.Text = ""

# Re-assert the autofit height the real deck ends up at (78.1pt ==
# 991870 EMU). Use a value that is infinitesimally above 78.1 so the
# point->EMU conversion rounds to the exact target instead of landing
# one EMU short because of binary floating point representation.
$sh.Height = 78.10001
